# Update countries & provincias Spain
# - Swap "Catar" / "Paises Bajos" order (row 23/24), refresh Catar's stats
# - Rotate "Malta" / "Republica de Africa Central" / "Etiopia" order (rows 130-132), refresh data
# - Bump the "Datos actualizados..." timestamp
# - Refresh several countries' case counts (rows 40, 57, 64, 66, 78)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 25 de Mayo de 2020 a las 13:35"

# --- Catar / Paises Bajos swap (rows 23-24) ---------------------------
$ws.Range("A23").Value = "Catar"
$ws.Range("B23").Value = 45465
$ws.Range("C23").Value = 1751
$ws.Range("D23").Value = 10363
$ws.Range("E23").Value = 35076
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 26

$ws.Range("A24").Value = "Paises Bajos"
$ws.Range("B24").Value = 45236
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("H24").Value = 5822

# --- Rumania refresh (row 40) ------------------------------------------
$ws.Range("E40").Value = 5456
$ws.Range("G40").Value = 12
$ws.Range("H40").Value = 1197

# --- Argelia refresh (row 57) -------------------------------------------
$ws.Range("D57").Value = 4578
$ws.Range("E57").Value = 3128

# --- Moldavia refresh (row 64) ------------------------------------------
$ws.Range("D64").Value = 3802
$ws.Range("E64").Value = 3035
$ws.Range("G64").Value = 6
$ws.Range("H64").Value = 256

# --- Finlandia refresh (row 66) -----------------------------------------
$ws.Range("D66").Value = 5100
$ws.Range("E66").Value = 1191
$ws.Range("G66").Value = 1
$ws.Range("H66").Value = 308

# --- Senegal refresh (row 78) -------------------------------------------
$ws.Range("B78").Value = 3130
$ws.Range("C78").Value = 83
$ws.Range("D78").Value = 1515
$ws.Range("E78").Value = 1580

# --- Etiopia / Malta / Republica de Africa Central rotation (rows 130-132) -
$ws.Range("A130").Value = "Etiopia"
$ws.Range("B130").Value = 655
$ws.Range("C130").Value = 73
$ws.Range("D130").Value = 159
$ws.Range("E130").Value = 491
$ws.Range("H130").Value = 5

$ws.Range("A131").Value = "Malta"
$ws.Range("B131").Value = 611
$ws.Range("C131").Value = 1
$ws.Range("D131").Value = 485
$ws.Range("E131").Value = 120
$ws.Range("H131").Value = 6

$ws.Range("A132").Value = "Republica de Africa Central"
$ws.Range("B132").Value = 604
$ws.Range("D132").Value = 22
$ws.Range("E132").Value = 581
$ws.Range("H132").Value = 1
